$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = "NSE:AGROPHOS"
$ws.Range("C2").Value = "NSE:AGI"
$ws.Range("D2").Value = $null
$ws.Range("E2").Value = "NSE:ADANIENT"
$ws.Range("F2").Value = $null

# Row 3
$ws.Range("B3").Value = "NSE:AUTOAXLES"
$ws.Range("C3").Value = "NSE:AURIONPRO"
$ws.Range("E3").Value = "NSE:ADANIGREEN"

# Row 4
$ws.Range("B4").Value = "NSE:GLOBUSSPR"
$ws.Range("C4").Value = "NSE:AURUM"
$ws.Range("E4").Value = "NSE:CAMS"

# Row 5
$ws.Range("B5").Value = "NSE:HEALTHY"
$ws.Range("C5").Value = "NSE:BVCL"
$ws.Range("E5").Value = "NSE:HFCL"

# Row 6
$ws.Range("B6").Value = "NSE:KALAMANDIR"
$ws.Range("C6").Value = "NSE:CERA"
$ws.Range("E6").Value = "NSE:IRCTC"

# Row 7
$ws.Range("B7").Value = "NSE:MON100"
$ws.Range("C7").Value = "NSE:CHOLAHLDNG"
$ws.Range("E7").Value = "NSE:OIL"

# Row 8
$ws.Range("B8").Value = "NSE:PHARMABEES"
$ws.Range("C8").Value = "NSE:CYIENTDLM"

# Row 9
$ws.Range("B9").Value = "NSE:PRIVISCL"
$ws.Range("C9").Value = "NSE:DOLLAR"

# Row 10
$ws.Range("B10").Value = $null
$ws.Range("C10").Value = "NSE:DPSCLTD"

# Row 11
$ws.Range("B11").Value = $null
$ws.Range("C11").Value = "NSE:GOLDIAM"

# Row 12
$ws.Range("B12").Value = $null
$ws.Range("C12").Value = "NSE:IDEAFORGE"

# Row 13
$ws.Range("B13").Value = $null
$ws.Range("C13").Value = "NSE:INDIACEM"

# Row 14
$ws.Range("B14").Value = $null
$ws.Range("C14").Value = "NSE:INDIAGLYCO"

# Row 15
$ws.Range("C15").Value = "NSE:INDSWFTLTD"

# Row 16
$ws.Range("C16").Value = "NSE:IONEXCHANG"

# Row 17
$ws.Range("C17").Value = "NSE:KOKUYOCMLN"

# Row 18
$ws.Range("C18").Value = "NSE:KTKBANK"

# Row 19
$ws.Range("C19").Value = "NSE:LOVABLE"

# Row 20
$ws.Range("C20").Value = "NSE:PITTIENG"

# Row 21
$ws.Range("C21").Value = "NSE:PRESTIGE"

# Row 22
$ws.Range("C22").Value = "NSE:QUICKHEAL"

# Row 23
$ws.Range("C23").Value = "NSE:RAJRATAN"

# Row 24
$ws.Range("C24").Value = "NSE:RAMCOCEM"

# Row 25
$ws.Range("C25").Value = "NSE:ROHLTD"

# Row 26
$ws.Range("C26").Value = "NSE:RTNPOWER"

# Remove the now-unused trailing rows 27-31 so the sheet's used range
# shrinks back down to A1:F26, matching the target dimension.
$ws.Range("A27:F31").Delete()
